$d = $word.ActiveDocument

$pairs = @(
    @("570×9=5130", "615×6=3690"),
    @("372×5=1860", "597×9=5373"),
    @("878×9=7902", "469×2=938"),
    @("461×7=3227", "256×3=768"),
    @("209×2=418", "128×2=256"),
    @("326×8=2608", "592×4=2368"),
    @("272×2=544", "790×4=3160"),
    @("941×7=6587", "951×4=3804"),
    @("703×6=4218", "388×9=3492"),
    @("645×8=5160", "512×7=3584"),
    @("455×7=3185", "783×3=2349"),
    @("331×7=2317", "195×3=585"),
    @("442×8=3536", "760×9=6840"),
    @("780×7=5460", "156×8=1248"),
    @("440×4=1760", "332×8=2656"),
    @("843×3=2529", "700×3=2100"),
    @("447×5=2235", "479×2=958"),
    @("878×6=5268", "219×7=1533"),
    @("628×4=2512", "720×5=3600"),
    @("181×2=362", "873×8=6984"),
    @("148×8=1184", "619×3=1857"),
    @("254×5=1270", "236×3=708"),
    @("422×4=1688", "718×6=4308"),
    @("155×8=1240", "648×5=3240"),
    @("196×2=392", "649×8=5192")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
